# RPA datasets push 2024-07-17
# The "노브랜드" (NoBrand) listing row (row 27) was removed from the dataset.
# Delete that entire row and shift the remaining rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(27).Delete()
